$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 -> becomes the "Blå slemspindling / Cortinarius salor" record (was row 6)
$ws.Range("A4").Value = 112363523
$ws.Range("B4").Value = 85252
$ws.Range("E4").Value = 3712
$ws.Range("F4").Value = "Blå slemspindling"
$ws.Range("G4").Value = "Cortinarius salor"
$ws.Range("H4").Value = "Fr."
$ws.Range("M4").Value = ""
$ws.Range("J4").Value = "fruktkroppar"
$ws.Range("Q4").Value = 721925
$ws.Range("R4").Value = 6397844

# Row 5 -> only the Taxonsorteringsordning (B) value updates
$ws.Range("B5").Value = 73772

# Row 6 -> becomes the "Blåmossa / Leucobryum glaucum" record (was row 7)
$ws.Range("A6").Value = 112363369
$ws.Range("B6").Value = 93553
$ws.Range("E6").Value = 2180
$ws.Range("F6").Value = "Blåmossa"
$ws.Range("G6").Value = "Leucobryum glaucum"
$ws.Range("H6").Value = "(Hedw.) Ångstr."
$ws.Range("J6").Value = "plantor/tuvor"
$ws.Range("Q6").Value = 721921
$ws.Range("R6").Value = 6397800

# Row 7 -> becomes the "Vågbandad barkbock / Semanotus undatus" record (was row 4)
$ws.Range("A7").Value = 112363550
$ws.Range("B7").Value = 5135
$ws.Range("E7").Value = 105930
$ws.Range("F7").Value = "Vågbandad barkbock"
$ws.Range("G7").Value = "Semanotus undatus"
$ws.Range("H7").Value = "(Linnaeus, 1758)"
$ws.Range("J7").Value = ""
$ws.Range("M7").Value = "färska gnagspår"
$ws.Range("Q7").Value = 721928
$ws.Range("R7").Value = 6397835
